$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

# Row 3 (CCGT set): add SetDesc (G) and And/Or markers (H/I)
$ws.Range("G3").Value = $ws.Range("F3").Text
$ws.Range("H3").Value = "And"
$ws.Range("I3").Value = "Or"

# Row 7 (OCGT set): add And/Or markers (H/I)
$ws.Range("H7").Value = "And"
$ws.Range("I7").Value = "Or"

# Row 17 (Nuclear set): add PSET_PN exclusion pattern (B) and And/Or markers (H/I)
$ws.Range("B17").Value = "-*SMR"
$ws.Range("H17").Value = "And"
$ws.Range("I17").Value = "Or"

# Extend PSET_PN wildcard patterns for CCGT and OCGT (gas turbine) rows
$ws.Range("B3").Value = $ws.Range("B3").Text + ",*GasCC*"
$ws.Range("B7").Value = $ws.Range("B7").Text + ",EN*CT*"
